$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "Weston"
$ws.Range("B10").Value = "Athena-Weston SD 29RJ"
$ws.Range("C10").Value = "Ann Vescio"
$ws.Range("D10").Value = "Laure Quraesma, Superintendent Athena-Weston Schools"
$ws.Range("E10").Value = "ann.vescio@athwest.k12.or.us"
$ws.Range("F10").Value = "laure.quaresma@athwest.k12.or.us"
$ws.Range("G10").Value = "541-566-3548"
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 5
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 3
$ws.Range("R10").Value = 4
$ws.Range("S10").Value = 2
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 2
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 2
$ws.Range("Y10").Value = 2
$ws.Range("Z10").Value = 2
$ws.Range("AA10").Value = "The Weston Middle School science department has worked diligently this year to align its science curriculum to the NGSS, to teach through an inquiry model,  and to regularly incorporate STEM activites. The Outdoor School program supported NGSS standards in the areas of PS3 Energy (Cooking Field Study); LS1 From Molecules to Organisms: Structures and Processes (Salmon Life Cycle Field Study and River Invertebrates Inquiry); LS2 Interactions, Energy and Dynamics Relationships in Ecosystems (Fire Impact Field Study, National Weather Service Field Study, Water Quality Field Study, Soil Field Study); LS3 Heredity: Inheritance and Variation of Traits (Salmon Life Cycle Field Study, River Invertebrates Inquiry); ESS2 Weather and Climate (National Weather Service Field Study, Forest Service Field Study); ESS3 Earth and Human Activity (Hike, Forest Service Field Study, Water Quality Field Study, Soil Field Study, River Invertebrates Inquiry); and ETS1 Engineering Design (Cooking Field Study).  Additionally, the water invertibrates, soil, and water quality presentations were inquiry-based lessons. The ""campfire cooking"" activity was actually a STEM activity which required students to design and then use a ""stove"" that they built to cook their lunches. These activities were based on the instructional model and on the NGSS curriculum that our students access daily in their science classes. "
$ws.Range("AB10").Value = "N/A"
$ws.Range("AC10").Value = "Weston Middle School is a small district located in a very rural area of Umatilla County in eastern Oregon. The school has a free/reduced lunch rate of 59%.  Outdoor School is an opportunity that few of our students would ever have through family resources or through district monies alone. The ability to provide Outdoor School to our students is one way that we can ensure that they have educational opportunities similar to students in wealthier areas of the state. "
$ws.Range("AD10").Value = "The field study sessions that we offered in the areas of salmon life cycle, water quality, forest ecosystems/forest fire impact, and soil conservation emphasized the interdendence of urban and rural areas. "
$ws.Range("AE10").Value = "X"
$ws.Range("AF10").Value = "X"
$ws.Range("AG10").Value = "NA"
$ws.Range("AH10").Value = "X"
$ws.Range("AI10").Value = "X"
$ws.Range("AJ10").Value = "X"
$ws.Range("AK10").Value = "X"
$ws.Range("AL10").Value = "Observation and Perspective (nature art/nature photography)"
$ws.Range("AM10").Value = "Students were required to turn in a completed Field Study Guide after their week of ODS. They debriefed in science classes about their learning at ODS. The two science teachers at Weston Middle School are attending a week-long science training this summer and are part of a science teaching cadre that meets three times annually. Ideas from these trainings will be incorporated into the plan for Outdoor School next school year. "
$ws.Range("AN10").Value = "Weston Middle School's Outdoor School program included the following partners: National Weather Service, Forest Service, Oregon State Land Laboratory, Department of Environmental Quality, CTUIR (tribe), two community volunteers, 20 parent volunteers, and a retired biology teacher volunteer. "
$ws.Range("AO10").Value = "N/A--Our ODS program was the last week of May this year. Testing had already been completed."
$ws.Range("AP10").Value = "The ODS program required students to assist one another and learn self-sufficiency in the First Aid session, the Water Invertibrates session, the cooking activity, and on the long hike. "
$ws.Range("AQ10").Value = "N/A--Our ODS program was the last week of May this year. There was no time to assess the impact on discipline and classroom management"
$ws.Range("AR10").Value = "Students were highly engaged in all facets of Outdoor School. There has been much positive feedback from both students and their parents about the success of the program. In addition, several of the presenters let the school know that the students were very well-behaved, engaged, and respectful. "
$ws.Range("AS10").Value = "Several activities promoted problem-solving and strategic thinking. The river invertebrate inquiry project required students to work together to collect invertebrates in a screen while wading in the river. The geocaching and first aid sessions required team work, problem-solving, and strategic thinking as well. Finally, the cooking session was an engineering project where students had to create a working ""campfire oven"" to cook their lunches. They had to use knowledge of energy and safety to create their ovens. "
$ws.Range("AT10").Value = "We included two Nature Art projects in our ODS program. One was to use items found on a hike to create frames for pictures. The other was to take creative photographs of nature using a variety of perspectives. "
$ws.Range("AU10").Value = "The students in our ODS program worked in groups throughout the week. Teachers worked with the teams prior to traveling to the Buck Creek facility to build group work skills through an orientation day. The sessions often required partner or team work to be successful. "
$ws.Range("AV10").Value = "N/A--Our ODS program was the last week of May this year. There was not time to determine whether or not increased enthusiasm for core subject areas transferred back into the classroom."
$ws.Range("AW10").Value = "The WMS science program has utilized an inquiry approach in the classroom throughout the year. In addition, students have participated in weekly Friday STEM labs all year. The students were able to apply their understanding of these processes in ODS in a variety of sessions including Water Invertebrates, Cooking, Water Quality, and Geocaching/Orienteering. In addition, students frequently utilized scientific vocabulary and background knowledge throughout the ODS program. "
$ws.Range("AX10").Value = "The cooking and first aid sessions in particular reinforced the application of scientific and civic processes in real-world situations. Students had to design a functional campfire ""oven"" with limited supplies in order to cook their stew for lunch. In the first aid sessions, students were given scenarios in which to apply their knowledge of first aid and rescue. "
$ws.Range("AY10").Value = "N/A"
$ws.Range("AZ10").Value = "Students kept notes in their ODS Field Study Guides as part of the program. They recorded notes and observations from all of their sessions following Notebooking protocols that they have utilized all year in their science and Friday STEM classes. "
$ws.Range("BA10").Value = "Students applied their knowledge of geography concepts in the Geocaching/Orienteering session which was run by their geography teacher. "
$ws.Range("BB10").Value = "All 6th graders but one who was on a family trip attended ODS. This included special education students as well as students of all abilities and learning styles. Support personnel attended to ensure that all programming was accessible to our special education students. "
$ws.Range("BC10").Value = "Students in the Athena-Weston School District come from a very rural area of Umatilla County. All students but one who was on a family trip attended ODS. The daytime-only format was important to our participation rate as we have a number of students and/or parents who indicated that they were uncomfortable with an overnight scenario. "
$ws.Range("BD10").Value = "N/A"
$ws.Range("BE10").Value = "All 6th grade students but one who was on a family trip attended ODS. This included our special education students. We sent adult support to ensure full participation from two autistic students and one Down's Syndrome student. "
$ws.Range("BF10").Value = "We had a student with Down's Syndrome participate in ODS. He has some physical limitations. Therefore we sent adult support to ensure that the program was accessible. "
$ws.Range("BG10").Value = "Weston Middle School's free and reduced lunch rate is 59%. We provided string packs, water bottles, binoculars, notebooks and writing utensils, t-shirts, reimbursable school lunches, water, and snacks at no charge to students to ensure that ODS was not a financial burden to any families. "
$ws.Range("BH10").Value = "21% of our 6th grade students identify as American Indian/Alaska Native. We included a presentation provided by the CTUIR tribal offices in our ODS program."
$ws.Range("BI10").Value = "1% of our 6th grade students identify as Asian."
$ws.Range("BJ10").Value = "N/A"
$ws.Range("BK10").Value = "1% of our 6th grade students identify as African American."
$ws.Range("BL10").Value = "11% of our 6th grade students identify as Hispanic/Latino."
$ws.Range("BM10").Value = "N/A"
$ws.Range("BN10").Value = "Transporting students to ODS daily during school hours versus an overnight program ensured that all 6th graders attended ODS. We will continue this format next year. "
$ws.Range("BO10").Value = "This was our first ODS program. Therefore we had no existing equipment for ODS. We purchased high quality radios so that the adult leaders could communicate as there was no reliable cell phone service in the Buck Creek area. This purchase would not have been feasible through our building budget. "
$ws.Range("BP10").Value = "Weston Middle School accessed local resources for its field study sessions. Those partners included the National Weather Service, CTUIR (tribal), OSU Land Lab scientists, DEQ Water Quality offices, and the Forest Service. In addition, a local retired biology teacher conducted a field study session. More than 20 parents volunteered throughout the week to chaperone, and four Weston Middle School teachers conducted field studies for ODS. We used the IMESD-owned Buck Creek facility at no charge for the ODS site."
